$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 108: VoiceDuration for "bonus_round" changes from 2 to 1.5
$ws.Range("C108").Value2 = 1.5

# Insert a new row at 112 for the new "division_cap" / "DIVISION" bonus-round
# cap entry (pushes the former rows 112-116 down to 113-117).
$ws.Rows.Item(112).Insert() | Out-Null
$ws.Range("A112").Value2 = "division_cap"
$ws.Range("B112").Value2 = "DIVISION"
$ws.Range("C112").Value2 = 1

# Append a new row at the end of the sheet (now row 118) for the
# link-disconnect-tutorial button text.
$ws.Range("A118").Value2 = "link_disconnect_tutorial"
$ws.Range("B118").Value2 = "Press this button to remove all the links."

# Match the author's final selection/scroll position in the sheet.
$ws.Range("B118").Select() | Out-Null
